$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("data")
$ws2 = $wb.Worksheets.Item("readme")

# Update the two shared-string text changes on the "readme" sheet
# ("to December 2024" -> "to December 2025", "Mar 90 - Dec 24" -> "Mar 90 - Dec 25")
$ws2.Range("B2").Value = "to December 2025"
$ws2.Range("A4").Value = "Mar 90 - Dec 25"

# Append the 12 new monthly VIX rows (Jan 2025 - Dec 2025) to the "data" sheet
$dates = @(45658, 45689, 45717, 45748, 45778, 45809, 45839, 45870, 45901, 45931, 45962, 45992)
$values = @(16.760000000000002, 16.97, 21.84, 31.97, 20.46, 18.399999999999999, 16.38, 15.75, 15.79, 18.09, 19.77, 15.55)

$startRow = 470
for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $startRow + $i
    $ws1.Cells.Item($row, 1).Value = $dates[$i]
    $ws1.Cells.Item($row, 2).Value = $values[$i]
}

# Selection / active-sheet / active-cell bookkeeping to mirror the saved view state
# (readme is selected first so that "data" ends up as the active/front sheet)
$ws2.Range("A5").Select()
$ws1.Range("B475").Select()
